$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the participant row for "Alyne Corona" (row 10), which shifts
# subsequent rows up and updates the used range/dimension accordingly.
$ws.Rows.Item(10).Delete()
